$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.714.81'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.852.86'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.44'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4268'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3591'
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07297'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8790'
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.79'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.855.01'
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.550'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.341'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06988'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.75'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008946'
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.630.60'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.002'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.053.52'
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.989'
$ws.Range("E25").Value = '  +4.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.57'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '120.12'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.251'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.880'
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08899'
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7601'
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.964'
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.516'
$ws.Range("E34").Value = '  -1.37%  '
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.110'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05424'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01930'
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.824'
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1672'
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5089'
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.618'
$ws.Range("E43").Value = '  -4.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.415'
$ws.Range("E44").Value = '  +1.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06526'
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '105.80'
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.33'
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.621'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.799'
$ws.Range("E51").Value = '  +1.21%  '
